$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of trade data (row 5) to the sheet, continuing the existing
# GILD trade data series.
$ws.Cells.Item(5, 1).Value = 9971.9500000000007
$ws.Cells.Item(5, 2).Value = 10013
$ws.Cells.Item(5, 3).Value = 80.11
$ws.Cells.Item(5, 4).Value = 79.78
$ws.Cells.Item(5, 5).Value = $false
$ws.Cells.Item(5, 6).Value = -0.41
$ws.Cells.Item(4, 7).Copy()
$ws.Cells.Item(5, 7).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(5, 7).Value = 42609.505266203705
$ws.Cells.Item(5, 8).Value = $false
